# Add voiceover filler notes to the slides that currently have no speaker
# notes (the notes placeholder is empty). This mirrors the commit
# "add voiceover filler for slides w/o voiceover".

$p = $ppt.ActivePresentation

$fillers = @{
    8  = "Exploratory Data Analysis"
    9  = "Preprocessing - NLP"
    10 = "Feature Engineering TBD"
    12 = "Generating Synthetic Data"
}

foreach ($idx in $fillers.Keys) {
    $slide = $p.Slides.Item($idx)
    $notesShape = $slide.NotesPage.Shapes.Item(2)
    $notesShape.TextFrame.TextRange.Text = $fillers[$idx]
}
